$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The "Inputs" mini-table (dates / price / cost-of-equity / growth) and
# the "Terminal value" mini-table currently live in column J (with a
# stray 4-column-wide block J:M reserved for them, and P:S left as
# unused filler cells). This edit moves both tables from column J back
# into column F ("Table in right column" -> really the values column
# sitting right next to the labels in column D/E), and drops the
# unused trailing P:S filler cells.
# ------------------------------------------------------------------

# 1) Move (cut/paste) the two 4-column-wide input blocks from J to F.
#    Using a 4-column source (J:M) - rather than just column J - makes
#    sure the destination picks up the correct per-cell formatting that
#    sat in the (empty) K/L/M companion cells too (e.g. the little
#    border accent that lived in L11 ends up on H11, matching the
#    original K/L/M -> G/H/I offset).
$ws.Range("J9:M13").Cut($ws.Range("F9"))
$ws.Range("J42:M49").Cut($ws.Range("F42"))

# 2) The cut above leaves the old formatting "echoed" behind in the
#    vacated J:M cells (this engine's Cut doesn't blank the source
#    formatting the way Excel's does). Reset those vacated cells back
#    to the plain/default look used by their neighbours.
foreach ($r in 9..13) {
    $ws.Range("G$r").Copy($ws.Range("J$r`:M$r"))
}
foreach ($r in 42..49) {
    $ws.Range("G$r").Copy($ws.Range("J$r`:M$r"))
}

# 3) Drop the unused filler cells in P:S for the affected rows entirely.
$ws.Range("P9:S13").Clear()
$ws.Range("P42:S49").Clear()

# 4) Re-point every formula that used to reference the J-column inputs
#    at their new F-column homes (the engine doesn't auto-repoint
#    formulas on cut/paste the way desktop Excel does).
$ws.Range("E17").Formula = "=F10"
$ws.Range("E18").Formula = "=F49"

$ws.Range("K25").Formula = "=`$F`$11"
$ws.Range("L25").Formula = "=`$F`$11"
$ws.Range("M25").Formula = "=`$F`$11"

$ws.Range("N34").Formula = "=(N33-`$F`$9)/365"
$ws.Range("O34").Formula = "=(O33-`$F`$9)/365"
$ws.Range("P34").Formula = "=(P33-`$F`$9)/365"
$ws.Range("Q34").Formula = "=(Q33-`$F`$9)/365"
$ws.Range("R34").Formula = "=(R33-`$F`$9)/365"

$ws.Range("N35").Formula = "=N31/(1+F12)^N34"
$ws.Range("O35").Formula = "=O31/(1+G12)^O34"
$ws.Range("P35").Formula = "=P31/(1+H12)^P34"
$ws.Range("Q35").Formula = "=Q31/(1+I12)^Q34"
$ws.Range("R35").Formula = "=R31/(1+J12)^R34"

$ws.Range("F42").Formula = "=R31"
$ws.Range("F43").Formula = "=F13"
$ws.Range("F44").Formula = "=F42*(1+F43)"
$ws.Range("F45").Formula = "=F44/(F12-F43)"
$ws.Range("F46").Formula = "=R34"
$ws.Range("F47").Formula = "=F45/(1+F46)^F12"
$ws.Range("F48").Formula = "=`$J`$37"
$ws.Range("F49").Formula = "=SUM(F47:F48)"

# 5) Sheet-view housekeeping: scroll position/selection moved, matching
#    the new focal point of the edit (J13 -> the growth-rate input cell
#    that used to sit at J13 before the move).
$ws.Activate()
$ws.Range("J13").Select()
